$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AmazonDataSheet")
$ws.Range("B2").Value = "Shoes"
